$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 690.25
$ws.Range("I12").Value = 913
$ws.Range("J12").Value = 22
$ws.Range("K12").Value = 913
$ws.Range("L12").Value = 22
$ws.Range("M12").Value = -743
$ws.Range("N12").Value = -362
$ws.Range("H62").Value = 7051.25
$ws.Range("I62").Value = 4205
$ws.Range("K62").Value = 4205
$ws.Range("M62").Value = -3581
$ws.Range("H65").Value = 7051.25
$ws.Range("I65").Value = 4205
$ws.Range("K65").Value = 21025
$ws.Range("M65").Value = -17905
$ws.Range("H80").Value = 285.7931
$ws.Range("I80").Value = 251.17647
$ws.Range("K80").Value = 753.52941
$ws.Range("M80").Value = 244.47059
$ws.Range("H83").Value = 285.7931
$ws.Range("I83").Value = 251.17647
$ws.Range("K83").Value = 2260.58823
$ws.Range("M83").Value = 2731.41177
$ws.Range("H132").Value = 2422.375
$ws.Range("I132").Value = 2339.8572
$ws.Range("K132").Value = 7019.571599999999
$ws.Range("M132").Value = -4489.571599999999
$ws.Range("H137").Value = 3130.389
$ws.Range("I137").Value = 1543.8889
$ws.Range("K137").Value = 4631.6667
$ws.Range("M137").Value = -2081.6667

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents() | Out-Null
$ws.Range("M2").ClearContents() | Out-Null
$ws.Range("N2").Value = 0
$ws.Range("H74").Value = 1330.2467
$ws.Range("I74").Value = 1038.3802
$ws.Range("K74").Value = 1038.3802
$ws.Range("M74").Value = -164.3802000000001
$ws.Range("H77").Value = 1330.2467
$ws.Range("I77").Value = 1038.3802
$ws.Range("K77").Value = 5191.901
$ws.Range("M77").Value = -823.9009999999998
$ws.Range("H88").Value = 597.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 597.25
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents() | Out-Null
$ws.Range("M88").Value = 597.25
$ws.Range("N88").Value = -1409.25
$ws.Range("H91").Value = 597.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 597.25
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents() | Out-Null
$ws.Range("M91").Value = 597.25
$ws.Range("N91").Value = -3405.25
$ws.Range("H110").Value = 9948.777
$ws.Range("I110").Value = 9953.166999999999
$ws.Range("K110").Value = 9953.166999999999
$ws.Range("M110").Value = -7908.166999999999
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").ClearContents() | Out-Null
$ws.Range("M116").ClearContents() | Out-Null
$ws.Range("N116").Value = 0

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents() | Out-Null
$ws.Range("M3").ClearContents() | Out-Null
$ws.Range("N3").Value = 0
$ws.Range("H22").Value = 667
$ws.Range("I22").Value = 667
$ws.Range("K22").Value = 667
$ws.Range("M22").Value = -494
$ws.Range("H132").Value = 112889.5
$ws.Range("J132").Value = 112889.5
$ws.Range("L132").Value = 112889.5
$ws.Range("N132").Value = -123009.5
$ws.Range("H134").Value = 2687.0667
$ws.Range("I134").Value = 2400.6667
$ws.Range("K134").Value = 7202.000100000001
$ws.Range("M134").Value = -4667.000100000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2432.1072
$ws.Range("I16").Value = 2343.0435
$ws.Range("K16").Value = 2343.0435
$ws.Range("M16").Value = -2056.0435
$ws.Range("H99").Value = 14685.6
$ws.Range("I99").Value = 11612.728
$ws.Range("J99").Value = 17100
$ws.Range("K99").Value = 11612.728
$ws.Range("L99").Value = 17100
$ws.Range("M99").Value = -10114.728
$ws.Range("N99").Value = -20096
$ws.Range("H113").Value = 2432.1072
$ws.Range("I113").Value = 2343.0435
$ws.Range("K113").Value = 2343.0435
$ws.Range("M113").Value = -173.0435000000002
$ws.Range("H122").Value = 2824.95
$ws.Range("I122").Value = 2685.75
$ws.Range("K122").Value = 8057.25
$ws.Range("M122").Value = -5607.25
$ws.Range("H126").Value = 14685.6
$ws.Range("I126").Value = 11612.728
$ws.Range("J126").Value = 17100
$ws.Range("K126").Value = 34838.18399999999
$ws.Range("L126").Value = 51300
$ws.Range("M126").Value = -32368.18399999999
$ws.Range("N126").Value = -56240

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").ClearContents() | Out-Null
$ws.Range("N105").Value = 0
$ws.Range("H129").Value = 2386.818
$ws.Range("I129").Value = 998
$ws.Range("K129").Value = 2994
$ws.Range("M129").Value = 2006
$ws.Range("H140").Value = 2137.5454
$ws.Range("I140").Value = 2137.5454
$ws.Range("K140").Value = 6412.6362
$ws.Range("M140").Value = -1232.6362

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3097.889
$ws.Range("J80").Value = 1996.2
$ws.Range("L80").Value = 1996.2
$ws.Range("N80").Value = -3992.2
$ws.Range("H83").Value = 3097.889
$ws.Range("J83").Value = 1996.2
$ws.Range("L83").Value = 9981
$ws.Range("N83").Value = -19965
$ws.Range("H113").Value = 4198.778
$ws.Range("I113").Value = 4671.5
$ws.Range("J113").Value = 4063.7144
$ws.Range("K113").Value = 4671.5
$ws.Range("L113").Value = 4063.7144
$ws.Range("M113").Value = -2501.5
$ws.Range("N113").Value = -8403.714400000001
$ws.Range("H126").Value = 4356.3335
$ws.Range("I126").Value = 3548.25
$ws.Range("K126").Value = 10644.75
$ws.Range("M126").Value = -8174.75
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents() | Out-Null
$ws.Range("N139").Value = 0

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1231.3077
$ws.Range("I16").Value = 1415
$ws.Range("J16").Value = 221
$ws.Range("K16").Value = 1415
$ws.Range("L16").Value = 221
$ws.Range("M16").Value = -1245
$ws.Range("N16").Value = -561

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2069.8462
$ws.Range("I81").Value = 1838
$ws.Range("J81").Value = 3345
$ws.Range("K81").Value = 3676
$ws.Range("L81").Value = 6690
$ws.Range("M81").Value = -2615
$ws.Range("N81").Value = -8812
$ws.Range("H84").Value = 2069.8462
$ws.Range("I84").Value = 1838
$ws.Range("J84").Value = 3345
$ws.Range("K84").Value = 18380
$ws.Range("L84").Value = 33450
$ws.Range("M84").Value = -13076
$ws.Range("N84").Value = -44058
$ws.Range("H122").Value = 1650.75
$ws.Range("I122").Value = 1034.3334
$ws.Range("K122").Value = 3103.0002
$ws.Range("M122").Value = -653.0001999999999
$ws.Range("H136").Value = 2018.303
$ws.Range("I136").Value = 1752.3103
$ws.Range("J136").Value = 3946.75
$ws.Range("K136").Value = 5256.9309
$ws.Range("M136").Value = -2706.9309
$ws.Range("N136").Value = -16940.25
